# Scheduled-runner update for the Golem Profits sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Refreshes cached market-price / leve-profit figures (columns H-N) for a
# handful of leve rows per sheet. Values below mirror the latest pull from
# the price-tracking job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

$ws.Range("H107").Value = 484.58823
$ws.Range("I107").Value = 487.46155
$ws.Range("K107").Value = 487.46155
$ws.Range("M107").Value = 1432.53845

$ws.Range("H112").Value = 5799
$ws.Range("J112").Value = 5998.75
$ws.Range("L112").Value = 17996.25
$ws.Range("N112").Value = -20212.25

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H138").Value = 1453
$ws.Range("I138").Value = 1453
$ws.Range("K138").Value = 4359
$ws.Range("M138").Value = 781

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()

$ws.Range("H26").Value = 4000
$ws.Range("J26").Value = 4333.3335
$ws.Range("L26").Value = 4333.3335
$ws.Range("N26").Value = -4993.3335

$ws.Range("H102").Value = 1889.8
$ws.Range("I102").Value = 1749.6666
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 1749.6666
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -127.6666
$ws.Range("N102").Value = -5344

$ws.Range("H124").Value = 47199.2
$ws.Range("J124").Value = 72000
$ws.Range("L124").Value = 72000
$ws.Range("N124").Value = -81820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 797.625
$ws.Range("I10").Value = 427.5
$ws.Range("K10").Value = 427.5
$ws.Range("M10").Value = -287.5

$ws.Range("H20").Value = 1440.1818
$ws.Range("I20").Value = 1142.75
$ws.Range("J20").Value = 2233.3333
$ws.Range("K20").Value = 1142.75
$ws.Range("L20").Value = 2233.3333
$ws.Range("M20").Value = -895.75
$ws.Range("N20").Value = -2727.3333

$ws.Range("H54").Value = 2370
$ws.Range("I54").Value = 2370
$ws.Range("K54").Value = 2370
$ws.Range("M54").Value = -1886

$ws.Range("H86").Value = 2398
$ws.Range("I86").Value = 2398
$ws.Range("K86").Value = 2398
$ws.Range("M86").Value = -1275

$ws.Range("H89").Value = 2398
$ws.Range("I89").Value = 2398
$ws.Range("K89").Value = 11990
$ws.Range("M89").Value = -6374

$ws.Range("H105").Value = 1150
$ws.Range("I105").Value = 1150
$ws.Range("K105").Value = 1150
$ws.Range("M105").Value = 597

$ws.Range("H134").Value = 1780.2667
$ws.Range("I134").Value = 1980.8334
$ws.Range("J134").Value = 978
$ws.Range("K134").Value = 5942.5002
$ws.Range("L134").Value = 2934
$ws.Range("M134").Value = -3407.5002
$ws.Range("N134").Value = -8004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1366.6666
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1366.6666
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = 1366.6666
$ws.Range("N3").Value = -1592.6666
$ws.Range("L3").ClearContents()

$ws.Range("H35").Value = 7416.5
$ws.Range("I35").Value = 7416.5
$ws.Range("K35").Value = 7416.5
$ws.Range("M35").Value = -7122.5

$ws.Range("H38").Value = 13724.714
$ws.Range("I38").Value = 2415
$ws.Range("J38").Value = 41999
$ws.Range("K38").Value = 2415
$ws.Range("L38").Value = 41999
$ws.Range("M38").Value = -2038
$ws.Range("N38").Value = -42753

$ws.Range("H46").Value = 13724.714
$ws.Range("I46").Value = 2415
$ws.Range("J46").Value = 41999
$ws.Range("K46").Value = 2415
$ws.Range("L46").Value = 41999
$ws.Range("M46").Value = -2204
$ws.Range("N46").Value = -42421

$ws.Range("H105").Value = 4412.25
$ws.Range("I105").Value = 4412.25
$ws.Range("K105").Value = 4412.25
$ws.Range("M105").Value = -2665.25

$ws.Range("H122").Value = 3111.2856
$ws.Range("I122").Value = 276.25
$ws.Range("K122").Value = 828.75
$ws.Range("M122").Value = 1621.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 40.666668
$ws.Range("I6").Value = 40.666668
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 122.000004
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -9.000004000000004
$ws.Range("M6").ClearContents()

$ws.Range("H7").Value = 101
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H50").Value = 500
$ws.Range("I50").Value = 500
$ws.Range("K50").Value = 1500
$ws.Range("M50").Value = -1019

$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5540

$ws.Range("H53").Value = 500
$ws.Range("I53").Value = 500
$ws.Range("K53").Value = 1500
$ws.Range("M53").Value = -1019

$ws.Range("H137").Value = 3281.5
$ws.Range("I137").Value = 1530
$ws.Range("J137").Value = 5033
$ws.Range("K137").Value = 4590
$ws.Range("L137").Value = 15099
$ws.Range("M137").Value = 510
$ws.Range("N137").Value = -25299

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.09090999999999
$ws.Range("I2").Value = 91.73333
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 91.73333
$ws.Range("L2").Value = 52
$ws.Range("M2").Value = 21.26667
$ws.Range("N2").Value = -278

$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1983.3334
$ws.Range("I22").Value = 1950
$ws.Range("K22").Value = 1950
$ws.Range("M22").Value = -1655

$ws.Range("H27").Value = 1983.3334
$ws.Range("I27").Value = 1950
$ws.Range("K27").Value = 1950
$ws.Range("M27").Value = -1843

$ws.Range("H32").Value = 7581.5
$ws.Range("I32").Value = 1775.3334
$ws.Range("K32").Value = 1775.3334
$ws.Range("M32").Value = -1458.3334

$ws.Range("H122").Value = 3198.111
$ws.Range("I122").Value = 2736.8
$ws.Range("J122").Value = 3774.75
$ws.Range("K122").Value = 8210.400000000001
$ws.Range("L122").Value = 11324.25
$ws.Range("M122").Value = -5760.400000000001
$ws.Range("N122").Value = -16224.25

$ws.Range("H124").Value = 22000
$ws.Range("J124").Value = 22000
$ws.Range("L124").Value = 22000
$ws.Range("N124").Value = -31820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1190
$ws.Range("I81").Value = 1190
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2380
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = -1319
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 1190
$ws.Range("I84").Value = 1190
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 11900
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = -6596
$ws.Range("M84").ClearContents()

$ws.Range("H132").Value = 1819.0834
$ws.Range("I132").Value = 1819.0834
$ws.Range("K132").Value = 5457.2502
$ws.Range("M132").Value = -2927.2502

$ws.Range("H136").Value = 1716.6666
$ws.Range("I136").Value = 1716.6666
$ws.Range("K136").Value = 5149.9998
$ws.Range("M136").Value = -2599.9998
